$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: sku / product_title / shopify_product_id / shopify_variant_id ---
$ws.Range("A2").Value2 = 722457911059
$ws.Range("B2").Value = "La Vie En Rose - 3 Item"
$ws.Range("C2").NumberFormat = "#######00000"
$ws.Range("C2").Value2 = 138427301906
$ws.Range("D2").NumberFormat = "#######00000"
$ws.Range("D2").Value2 = 1340705931282

# --- Row 3: sku / product_title / shopify_product_id / shopify_variant_id ---
$ws.Range("A3").Value2 = 722457572946
$ws.Range("B3").Value = "La Vie En Rose - 5 Item"
$ws.Range("C3").NumberFormat = "#######00000"
$ws.Range("C3").Value2 = 138427203602
$ws.Range("D3").NumberFormat = "#######00000"
$ws.Range("D3").Value2 = 1340704751634

# --- New column E: product_collection (mirrors product_title) ---
$ws.Range("E1").Value = "product_collection"
$ws.Range("E2").Value = "La Vie En Rose - 3 Item"
$ws.Range("E3").Value = "La Vie En Rose - 5 Item"

# --- Match author's final selection state ---
$ws.Range("E3").Select() | Out-Null
